$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05340925434738789
$ws.Range("C2").Value = 0.9984293249484364
$ws.Range("D2").Value = 0.1687063053803154
$ws.Range("G2").Value = 0.1248244242667473
$ws.Range("H2").Value = 0.99
